$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: paragraph "возмездного оказания услуг № {{contract_number_contents}}
# от {{contract_start_date_contents}} года" -> collapse the merge-field runs into
# a single literal run, and drop the stray <w:lang w:val="en-US"/> from the
# paragraph-mark run properties.
# ------------------------------------------------------------------
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -match "contract_number_contents") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    $p1 = $d.Paragraphs($targetIndex)
    $xml1 = "<w:p><w:pPr><w:jc w:val=`"center`"/><w:rPr><w:sz w:val=`"22`"/><w:szCs w:val=`"22`"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val=`"22`"/><w:szCs w:val=`"22`"/></w:rPr><w:t>возмездного оказания услуг № ДП4000584 от 23.06.2020 года</w:t></w:r></w:p>"
    $p1.Range.InsertXML($xml1)

    # --------------------------------------------------------------
    # Change 2: the next (empty) paragraph keeps its <w:b/> paragraph-mark
    # run property but loses its stray <w:lang w:val="en-US"/>.
    # --------------------------------------------------------------
    $p2 = $d.Paragraphs($targetIndex + 1)
    $xml2 = "<w:p><w:pPr><w:spacing w:line=`"276`" w:lineRule=`"auto`"/><w:rPr><w:b/></w:rPr></w:pPr></w:p>"
    $p2.Range.InsertXML($xml2)

    # --------------------------------------------------------------
    # Change 3: the following paragraph ("Настоящим письмом ...") loses its
    # <w:jc w:val="both"/> paragraph alignment (reverts to the Word default,
    # i.e. the attribute disappears entirely instead of becoming "left").
    # --------------------------------------------------------------
    $p3 = $d.Paragraphs($targetIndex + 2)
    $p3.Format.Alignment = 0
}

# ------------------------------------------------------------------
# Change 4: "Просим обеспечить доступ в ресторан." gains a merge field for
# the access date: "Просим обеспечить доступ в ресторан {{access_date}}."
# ------------------------------------------------------------------
$accessIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -match "Просим обеспечить доступ в ресторан") {
        $accessIndex = $i
        break
    }
}

if ($accessIndex -gt 0) {
    $p4 = $d.Paragraphs($accessIndex)
    $xml4 = "<w:p>" +
              "<w:pPr>" +
                "<w:jc w:val=`"both`"/>" +
                "<w:rPr><w:bCs/><w:sz w:val=`"22`"/><w:szCs w:val=`"22`"/></w:rPr>" +
              "</w:pPr>" +
              "<w:r>" +
                "<w:rPr><w:bCs/><w:sz w:val=`"22`"/><w:szCs w:val=`"22`"/></w:rPr>" +
                "<w:t>Просим обеспечить доступ в ресторан</w:t>" +
              "</w:r>" +
              "<w:r>" +
                "<w:rPr><w:bCs/><w:sz w:val=`"22`"/><w:szCs w:val=`"22`"/></w:rPr>" +
                "<w:t xml:space=`"preserve`"> </w:t>" +
              "</w:r>" +
              "<w:r>" +
                "<w:rPr><w:color w:val=`"000000`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"22`"/></w:rPr>" +
                "<w:t>{{</w:t>" +
              "</w:r>" +
              "<w:r>" +
                "<w:rPr><w:color w:val=`"000000`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"22`"/><w:lang w:val=`"en-US`"/></w:rPr>" +
                "<w:t>access</w:t>" +
              "</w:r>" +
              "<w:r>" +
                "<w:rPr><w:color w:val=`"000000`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"22`"/></w:rPr>" +
                "<w:t>_</w:t>" +
              "</w:r>" +
              "<w:r>" +
                "<w:rPr><w:color w:val=`"000000`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"22`"/><w:lang w:val=`"en-US`"/></w:rPr>" +
                "<w:t>date</w:t>" +
              "</w:r>" +
              "<w:r>" +
                "<w:rPr><w:color w:val=`"000000`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"22`"/></w:rPr>" +
                "<w:t>}</w:t>" +
              "</w:r>" +
              "<w:r>" +
                "<w:rPr><w:color w:val=`"000000`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"22`"/></w:rPr>" +
                "<w:t>}</w:t>" +
              "</w:r>" +
              "<w:r>" +
                "<w:rPr><w:bCs/><w:sz w:val=`"22`"/><w:szCs w:val=`"22`"/></w:rPr>" +
                "<w:t>.</w:t>" +
              "</w:r>" +
            "</w:p>"
    $p4.Range.InsertXML($xml4)
}
